# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the latest generated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 70
$ws1.Range("F4").Value  = 1500
$ws1.Range("F5").Value  = 578
$ws1.Range("F6").Value  = 1063
$ws1.Range("F7").Value  = 11045
$ws1.Range("F8").Value  = 11045
$ws1.Range("F11").Value = 319
$ws1.Range("F14").Value = 12221
$ws1.Range("F15").Value = 12752
$ws1.Range("F16").Value = 32
$ws1.Range("F17").Value = 124
$ws1.Range("F20").Value = 78
$ws1.Range("F22").Value = 37

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 70
$ws4.Range("F5").Value  = 1500
$ws4.Range("F6").Value  = 578
$ws4.Range("F7").Value  = 1063
$ws4.Range("F8").Value  = 11045
$ws4.Range("F9").Value  = 11046
$ws4.Range("F12").Value = 319
$ws4.Range("F15").Value = 12221
$ws4.Range("F16").Value = 12752
$ws4.Range("F17").Value = 32
$ws4.Range("F18").Value = 124
$ws4.Range("F21").Value = 78
$ws4.Range("F23").Value = 37
